{"js": "// Applies the copy edits described in the commit:\n//  - \"Through out\" -> \"Throughout\"\n//  - removes \"in this class \" from the Object-Oriented Design sentence\n//  - adds a comma after \"During this time\" (first/MVC card occurrence only)\n//  - \"the majority of ... three course\" -> \"most of ... three courses\"\n//  - \"photo from a cellphone\" -> \"photos from a cellphone\"\n//  - \"some day\" -> \"someday\"\n//\n// Each replacement searches for a long, unique phrase so we never touch the\n// (very similarly worded) \"During this time,\" sentence that already existed\n// in the Entity Framework Core card of the source document.\n\nasync function replaceOnce(body, searchText, newText, options) {\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \"Through out\" -> \"Throughout\"\nawait replaceOnce(\n  body,\n  \"Through out my entire time in school, we spent a great deal of time\",\n  \"Throughout my entire time in school, we spent a great deal of time\"\n);\n\n// 2) remove \"in this class \" before \"we were tasked with creating our own program\"\nawait replaceOnce(\n  body,\n  \" was called Object-Oriented Design and in this class we were tasked with creating our own program from scratch using\",\n  \" was called Object-Oriented Design and we were tasked with creating our own program from scratch using\"\n);\n\n// 3) add comma: \"During this time\" -> \"During this time,\" (first occurrence only \u2013\n//    the MVC/ASP.Net \"Read more\" card; the Entity Framework Core card already has the comma)\nawait replaceOnce(\n  body,\n  \"During this time we also worked on a group project of our own choosing and were tasked with doing the same thing, but as a group without much guidance.\",\n  \"During this time, we also worked on a group project of our own choosing and were tasked with doing the same thing, but as a group without much guidance.\"\n);\n\n// 4) \"the majority of ... three course\" -> \"most of ... three courses\"\nawait replaceOnce(\n  body,\n  \"Although we were taught the majority of our knowledge in the three course mentioned above\",\n  \"Although we were taught most of our knowledge in the three courses mentioned above\"\n);\n\n// 5) \"photo from a cellphone\" -> \"photos from a cellphone\"\nawait replaceOnce(\n  body,\n  \"didn\u2019t come across as bland as photo from a cellphone or the old digital cameras I grew up with.\",\n  \"didn\u2019t come across as bland as photos from a cellphone or the old digital cameras I grew up with.\"\n);\n\n// 6) \"some day\" -> \"someday\"\nawait replaceOnce(\n  body,\n  \"Maybe some day I\u2019ll try to turn this hobby into more of a passion\",\n  \"Maybe someday I\u2019ll try to turn this hobby into more of a passion\"\n);\n", "ps1": "# Applies the copy edits described in the commit:\n#  - \"Through out\" -> \"Throughout\"\n#  - removes \"in this class \" from the Object-Oriented Design sentence\n#  - adds a comma after \"During this time\" (first/MVC card occurrence only)\n#  - \"the majority of ... three course\" -> \"most of ... three courses\"\n#  - \"photo from a cellphone\" -> \"photos from a cellphone\"\n#  - \"some day\" -> \"someday\"\n#\n# Each replacement searches on a long, unique phrase so we never touch the\n# (very similarly worded) \"During this time,\" sentence that already existed\n# in the Entity Framework Core card of the source document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Search text not found: $findText\"\n    }\n}\n\n# 1) \"Through out\" -> \"Throughout\"\nReplace-Once \"Through out my entire time in school, we spent a great deal of time\" \"Throughout my entire time in school, we spent a great deal of time\"\n\n# 2) remove \"in this class \" before \"we were tasked with creating our own program\"\nReplace-Once \" was called Object-Oriented Design and in this class we were tasked with creating our own program from scratch using\" \" was called Object-Oriented Design and we were tasked with creating our own program from scratch using\"\n\n# 3) add comma: \"During this time\" -> \"During this time,\" (first occurrence only -\n#    the MVC/ASP.Net \"Read more\" card; the Entity Framework Core card already has the comma)\nReplace-Once \"During this time we also worked on a group project of our own choosing and were tasked with doing the same thing, but as a group without much guidance.\" \"During this time, we also worked on a group project of our own choosing and were tasked with doing the same thing, but as a group without much guidance.\"\n\n# 4) \"the majority of ... three course\" -> \"most of ... three courses\"\nReplace-Once \"Although we were taught the majority of our knowledge in the three course mentioned above\" \"Although we were taught most of our knowledge in the three courses mentioned above\"\n\n# 5) \"photo from a cellphone\" -> \"photos from a cellphone\"\nReplace-Once \"didn\u2019t come across as bland as photo from a cellphone or the old digital cameras I grew up with.\" \"didn\u2019t come across as bland as photos from a cellphone or the old digital cameras I grew up with.\"\n\n# 6) \"some day\" -> \"someday\"\nReplace-Once \"Maybe some day I\u2019ll try to turn this hobby into more of a passion\" \"Maybe someday I\u2019ll try to turn this hobby into more of a passion\"\n"}
